# Update "PERIOD TO EXPIRE" (column H) and "LAST UPDATE" (column I)
# on the Training Dashboard sheet to reflect progress as of 04-Nov-2025.
#
# Column I holds the date as literal text (e.g. "04-Nov-2025"), not a real
# date value. Assigning a date-looking string straight to .Value would make
# Excel auto-convert it into a date serial number (and change the cell's
# number format/style). To keep it as plain text with the original style,
# we first write it as a formula that evaluates to the literal string, then
# use Copy + PasteSpecial(values) to collapse the formula down to a static
# text value while preserving the existing cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$xlPasteValues = -4163

# Row 3
$ws.Range("H3").Value = 345
$c = $ws.Range("I3")
$c.Formula = "=""04-Nov-2025"""
$c.Copy()
$c.PasteSpecial($xlPasteValues)

# Row 4
$ws.Range("H4").Value = 364
$c = $ws.Range("I4")
$c.Formula = "=""04-Nov-2025"""
$c.Copy()
$c.PasteSpecial($xlPasteValues)

# Row 5
$ws.Range("H5").Value = 364
$c = $ws.Range("I5")
$c.Formula = "=""04-Nov-2025"""
$c.Copy()
$c.PasteSpecial($xlPasteValues)

# Row 6
$ws.Range("H6").Value = 604
$c = $ws.Range("I6")
$c.Formula = "=""04-Nov-2025"""
$c.Copy()
$c.PasteSpecial($xlPasteValues)

$excel.CutCopyMode = $false
